# Updated symbol list (Price / Volume(1h) columns) with refreshed coinranking.com
# quotes. Values are written with a leading apostrophe so Excel stores them as
# text (matching the original inlineStr cells) instead of silently coercing
# numeric-looking strings to Numbers, which would drop significant trailing
# zeros (e.g. D26: 0.0001302 -> 0.0001300).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.88"
$ws.Range("E2").Value = "'-0.20%"
$ws.Range("D3").Value = "'41.04"
$ws.Range("E3").Value = "'0.69%"
$ws.Range("D4").Value = "'5.212"
$ws.Range("E4").Value = "'1.85%"
$ws.Range("D5").Value = "'0.07672"
$ws.Range("E5").Value = "'0.61%"
$ws.Range("D6").Value = "'1.633"
$ws.Range("E6").Value = "'1.35%"
$ws.Range("D7").Value = "'0.9149"
$ws.Range("E7").Value = "'1.20%"
$ws.Range("D9").Value = "'0.1211"
$ws.Range("E9").Value = "'9.42%"
$ws.Range("E10").Value = "'3.16%"
$ws.Range("D11").Value = "'0.09152"
$ws.Range("E11").Value = "'-0.02%"
$ws.Range("D12").Value = "'0.04158"
$ws.Range("E12").Value = "'-0.57%"
$ws.Range("E13").Value = "'-0.29%"
$ws.Range("D14").Value = "'0.001261"
$ws.Range("E14").Value = "'0.16%"
$ws.Range("D15").Value = "'0.005769"
$ws.Range("E15").Value = "'-1.84%"
$ws.Range("D17").Value = "'3.341"
$ws.Range("D18").Value = "'4.308"
$ws.Range("E18").Value = "'1.31%"
$ws.Range("D19").Value = "'0.3336"
$ws.Range("E19").Value = "'1.32%"
$ws.Range("D20").Value = "'7.394"
$ws.Range("E20").Value = "'12.75%"
$ws.Range("D23").Value = "'0.04018"
$ws.Range("E23").Value = "'-1.00%"
$ws.Range("D24").Value = "'0.001258"
$ws.Range("E24").Value = "'2.27%"
$ws.Range("D25").Value = "'0.004380"
$ws.Range("E25").Value = "'6.46%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'0.02%"
$ws.Range("D38").Value = "'0.02492"
$ws.Range("D39").Value = "'0.05329"
$ws.Range("E39").Value = "'2.71%"
$ws.Range("D40").Value = "'0.007841"
$ws.Range("E40").Value = "'0.80%"
$ws.Range("D41").Value = "'0.1312"
$ws.Range("E41").Value = "'0.49%"
$ws.Range("D42").Value = "'0.006509"
$ws.Range("E42").Value = "'-3.50%"
$ws.Range("D43").Value = "'0.001910"
$ws.Range("E43").Value = "'-2.03%"
$ws.Range("D44").Value = "'0.008249"
$ws.Range("E44").Value = "'-6.08%"
$ws.Range("D45").Value = "'0.3338"
$ws.Range("E45").Value = "'0.18%"
$ws.Range("D46").Value = "'0.00006705"
$ws.Range("E46").Value = "'-3.63%"
$ws.Range("E47").Value = "'0.09%"
$ws.Range("D48").Value = "'0.3526"
$ws.Range("E48").Value = "'1,049.50%"
$ws.Range("D49").Value = "'0.003105"
$ws.Range("E49").Value = "'-26.13%"
$ws.Range("E50").Value = "'0.09%"
$ws.Range("E51").Value = "'0.09%"
